$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.214.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.25%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.273.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'583.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.07%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'185.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.06%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.58%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.47%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.844.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.15%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.198.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.49%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.268.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'5.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.22%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'414.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.08%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.63%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.09%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.34%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.509"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.53%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.01%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.13%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'22.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.85%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.37%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E34").Value = "'  -1.87%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'164.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.49%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.92%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +4.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.55%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.659.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'40.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.98%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0680"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'338.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'24.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0275"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.83%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.974"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.43%  "
$ws.Range("E51").Style = "Normal"
